# Update "想去人数" (want-to-go count) figures in column F across sheets.
# This mirrors a data refresh (gh-pages output regeneration) where several
# event rows had their interested-attendee counts incremented.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 146
$ws1.Range("F5").Value  = 3324
$ws1.Range("F6").Value  = 1050
$ws1.Range("F7").Value  = 2204
$ws1.Range("F8").Value  = 2117
$ws1.Range("F10").Value = 606
$ws1.Range("F11").Value = 21
$ws1.Range("F16").Value = 94
$ws1.Range("F17").Value = 210
$ws1.Range("F19").Value = 632
$ws1.Range("F20").Value = 722
$ws1.Range("F21").Value = 603
$ws1.Range("F22").Value = 12243
$ws1.Range("F23").Value = 12286
$ws1.Range("F25").Value = 702
$ws1.Range("F27").Value = 33
$ws1.Range("F28").Value = 20
$ws1.Range("F29").Value = 366
$ws1.Range("F30").Value = 1920
$ws1.Range("F31").Value = 4
$ws1.Range("F32").Value = 197
$ws1.Range("F33").Value = 582

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 26

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 146
$ws4.Range("F6").Value  = 3324
$ws4.Range("F7").Value  = 1050
$ws4.Range("F8").Value  = 2204
$ws4.Range("F9").Value  = 2117
$ws4.Range("F11").Value = 606
$ws4.Range("F12").Value = 21
$ws4.Range("F19").Value = 94
$ws4.Range("F21").Value = 210
$ws4.Range("F23").Value = 632
$ws4.Range("F24").Value = 722
$ws4.Range("F25").Value = 603
$ws4.Range("F26").Value = 12243
$ws4.Range("F27").Value = 12286
$ws4.Range("F29").Value = 702
$ws4.Range("F31").Value = 33
$ws4.Range("F32").Value = 20
$ws4.Range("F33").Value = 366
$ws4.Range("F34").Value = 1920
$ws4.Range("F36").Value = 4
$ws4.Range("F38").Value = 197
$ws4.Range("F39").Value = 582
$ws4.Range("F40").Value = 26
